$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 964, pushing existing row 964 (and below)
# down to row 965, etc. This mirrors the diff: a new price record is
# inserted in the middle of the weekly dataset, shifting every
# subsequent row index by +1 (old row 1040 becomes row 1041).
$ws.Rows.Item(964).Insert()

# Populate the newly inserted row 964 with the new record's data.
$ws.Cells.Item(964, 1).Value = 3
$ws.Cells.Item(964, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(964, 3).Value = "Coquimbo"
$ws.Cells.Item(964, 4).Value = 45021
$ws.Cells.Item(964, 5).Value = 5
$ws.Cells.Item(964, 6).Value = 100112024
$ws.Cells.Item(964, 7).Value = "Choclo"
$ws.Cells.Item(964, 8).Value = "Choclero"
$ws.Cells.Item(964, 9).Value = "Primera"
$ws.Cells.Item(964, 10).Value = 13500
$ws.Cells.Item(964, 11).Value = 350
$ws.Cells.Item(964, 12).Value = 350
$ws.Cells.Item(964, 13).Value = 350
$ws.Cells.Item(964, 14).Value = "`$/unidad"
$ws.Cells.Item(964, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(964, 16).Value = 350
$ws.Cells.Item(964, 17).Value = 1
$ws.Cells.Item(964, 18).Value = "Hortaliza"
